$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Global font-name fix: TimesNewToman -> Times New Roman
#    Done FIRST (format-only find/replace) so that any runs split later by
#    this script inherit the corrected font automatically.
# ---------------------------------------------------------------------------
$fontFind = $d.Content.Find
$fontFind.ClearFormatting()
$fontFind.Font.Name = "TimesNewToman"
$fontFind.Replacement.ClearFormatting()
$fontFind.Replacement.Font.Name = "Times New Roman"
$fontFind.Execute("", $false, $false, $false, $false, $false, $true, 1, $true, "", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Title
# ---------------------------------------------------------------------------
$f = $d.Content.Find
$f.ClearFormatting()
$f.Execute("Exploring the Enigmatic Realm of Dreams", $false, $false, $false, $false, $false, $true, 1, $false, "The Marvelous World of Biology: Exploring the Symphony of Life", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Author name: "Emma Williams" -> "Dr. Emily Carter" split across 3 runs
#    ("Dr", ".", " Emily Carter") to mirror the target markup.
# ---------------------------------------------------------------------------
$f = $d.Content.Find
$f.ClearFormatting()
$f.Execute("Emma Williams", $false, $false, $false, $false, $false, $true, 1, $false, "Dr. Emily Carter", 2) | Out-Null

$namePara = $d.Paragraphs.Item(2).Range
$dotStart = $namePara.Start + 2
$dotRange = $d.Range($dotStart, $dotStart + 1)
# Toggling a direct character attribute forces Word to split the run at this
# sub-string boundary; flipping it back off keeps the (now 3-way) run split
# without leaving any visible formatting difference.
$dotRange.Bold = 1
$dotRange.Bold = 0

# ---------------------------------------------------------------------------
# 4) Email address parts
# ---------------------------------------------------------------------------
$f = $d.Content.Find
$f.ClearFormatting()
$f.Execute("emma", $true, $true, $false, $false, $false, $true, 1, $false, "carter", 2) | Out-Null

$f = $d.Content.Find
$f.ClearFormatting()
$f.Execute("williams@berkeley", $false, $false, $false, $false, $false, $true, 1, $false, "emily725@edumail", 2) | Out-Null

$f = $d.Content.Find
$f.ClearFormatting()
$f.Execute("edu", $true, $true, $false, $false, $false, $true, 1, $false, "org", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Body paragraph (dreams -> biology)
# ---------------------------------------------------------------------------
$f = $d.Content.Find
$f.ClearFormatting()
$f.Execute("As humans, we spend a significant portion of our lives in the ethereal realm of dreams, navigating landscapes both familiar and fantastical", $false, $false, $false, $false, $false, $true, 1, $false, "Biology, the study of life, unveils a captivating world teeming with intricate processes, marvelous adaptations, and awe-inspiring diversity", 2) | Out-Null

$f = $d.Content.Find
$f.ClearFormatting()
$f.Execute(" Dreams have captivated the imaginations of philosophers, artists, and scientists alike throughout history, inspiring countless works of art, literature, and scientific inquiry", $false, $false, $false, $false, $false, $true, 1, $false, " This scientific discipline invites us to delve into the secrets held within living organisms, from the smallest bacteria to the colossal blue whale", 2) | Out-Null

$f = $d.Content.Find
$f.ClearFormatting()
$f.Execute(" From Sigmund Freud's psychoanalytic theories to the latest findings in neuroscience, the quest to understand the nature and significance of dreams continues to be an enthralling intellectual pursuit", $false, $false, $false, $false, $false, $true, 1, $false, " In this journey, we uncover the symphony of life, the delicate balance of ecosystems, and the extraordinary resilience of life on Earth", 2) | Out-Null

$f = $d.Content.Find
$f.ClearFormatting()
$f.Execute("Delving into the labyrinthine corridors of the dreaming mind, we encounter a world governed by its own unique logic, where time and space are fluid, and the boundaries between reality and imagination blur", $false, $false, $false, $false, $false, $true, 1, $false, "Biology unravels the genetic code, the blueprint that guides the development, structure, and function of all living beings", 2) | Out-Null

$f = $d.Content.Find
$f.ClearFormatting()
$f.Execute(" Dreams transport us to strange and wonderful places, introduce us to unforgettable characters, and often leave us with lingering emotions and insights upon waking", $false, $false, $false, $false, $false, $true, 1, $false, " We marvel at the intricate dance of molecules, the conversion of energy, and the remarkable symphony of chemical reactions that orchestrate the processes of life", 2) | Out-Null

$f = $d.Content.Find
$f.ClearFormatting()
$f.Execute(" These nocturnal journeys have long been a source of fascination, a window into the hidden recesses of our subconscious minds", $false, $false, $false, $false, $false, $true, 1, $false, " Through microscopy, we peer into the hidden world of cells, uncovering their diverse structures and specialized functions, revealing a universe within each minuscule entity", 2) | Out-Null

$f = $d.Content.Find
$f.ClearFormatting()
$f.Execute("Dreams have been the subject of scientific scrutiny for centuries, with researchers employing a variety of methods to unravel their mysteries", $false, $false, $false, $false, $false, $true, 1, $false, "Furthermore, biology illuminates the interdependence of organisms, the delicate web of life that connects all ecosystems", 2) | Out-Null

$f = $d.Content.Find
$f.ClearFormatting()
$f.Execute(" From detailed dream diaries to advanced neuroimaging techniques, scientists are gradually piecing together the intricate puzzle of why we dream", $false, $false, $false, $false, $false, $true, 1, $false, " We discover the intricate interactions between species, the delicate balance of predators and prey, and the vital role of biodiversity in maintaining the health of our planet", 2) | Out-Null

$f = $d.Content.Find
$f.ClearFormatting()
$f.Execute(" While the precise functions of dreams remain elusive, research suggests that they play a role in memory consolidation, emotional regulation, and creative problem-solving", $false, $false, $false, $false, $false, $true, 1, $false, " As we explore the diversity of life, we gain a deeper understanding of our own place in the intricate tapestry of life on Earth", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6) Summary paragraph
# ---------------------------------------------------------------------------
$f = $d.Content.Find
$f.ClearFormatting()
$f.Execute("Dreams, those enigmatic nocturnal experiences, have intrigued humanity for millennia", $false, $false, $false, $false, $false, $true, 1, $false, "Biology, the study of life, captivates us with its exploration of the intricate symphony of living organisms, the genetic code that guides their existence, the diversity of ecosystems, and the remarkable resilience of life", 2) | Out-Null

$f = $d.Content.Find
$f.ClearFormatting()
$f.Execute(" From artistic and literary musings to scientific investigations, the study of dreams has shed light on the hidden workings of our minds", $false, $false, $false, $false, $false, $true, 1, $false, " It unveils the delicate dance of molecules, the intricate structures of cells, and the interdependence of organisms", 2) | Out-Null

$f = $d.Content.Find
$f.ClearFormatting()
$f.Execute(" While the precise purpose of dreams remains a subject of ongoing research, their role in memory consolidation, emotional regulation, and creative thinking is increasingly recognized", $false, $false, $false, $false, $false, $true, 1, $false, " Biology reveals the wonders of adaptation, the marvels of diversity, and the extraordinary resilience of life, expanding our understanding of the interconnectedness of all living beings and inspiring us to appreciate the beauty and fragility of our planet", 2) | Out-Null

# Drop the trailing sentence (". As we continue ... unconscious") while
# keeping the paragraph's very last, standalone "." run intact.
$marker = $d.Content.Find
$marker.ClearFormatting()
$marker.Text = "planet"
$marker.Forward = $true
$marker.Wrap = 0
$marker.Execute() | Out-Null
$markerRange = $marker.Parent

$summaryPara = $d.Paragraphs.Item(7).Range
$deleteStart = $markerRange.End
$deleteEnd = $summaryPara.End - 2
$trailing = $d.Range($deleteStart, $deleteEnd)
$trailing.Delete()

# Re-split the final "." back into its own run (formatting-identical runs
# otherwise get coalesced by the text-delete above).
$summaryPara = $d.Paragraphs.Item(7).Range
$lastDotStart = $summaryPara.End - 2
$lastDotRange = $d.Range($lastDotStart, $lastDotStart + 1)
$lastDotRange.Bold = 1
$lastDotRange.Bold = 0

# ---------------------------------------------------------------------------
# 7) Trailing empty paragraph added at the very end of the document
# ---------------------------------------------------------------------------
$d.Paragraphs.Last.Range.InsertParagraphAfter()

Write-Output "done"
